$wb = $excel.ActiveWorkbook

# Duplicate the "Turkey" sheet (same layout/column widths as the new Greece sheet)
# and place the copy after the last sheet (Croatia).
$turkey = $wb.Worksheets.Item("Turkey")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$turkey.Copy($null, $lastSheet)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Greece"

# Update the market name / reference values for Greece
$ws.Range("B2").Value = "Greece Market"
$ws.Range("B4").Value = "NGC-4119/T3208/3192"

# Croatia is no longer the active tab; its selection becomes the full sheet select
$croatia = $wb.Worksheets.Item("Croatia")
[void]$croatia.Cells.Select()

# Select / activate B4 on the new sheet last so it becomes the active tab
$ws.Activate()
[void]$ws.Range("B4").Select()
